$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M13").Value = 4881.28

# Sheet: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 4881.28
$ws2.Range("F23").Value = 7205.26

# Sheet: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 6842.21
$ws3.Range("E12").Value = 37575.79
$ws3.Range("F12").Value = 0.1540413796208744

$ws3.Range("D14").Value = 7205.26
$ws3.Range("E14").Value = 48194.21101170094
$ws3.Range("F14").Value = 0.1300600866473648
